$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (GLN 7601000157638 / "à Wengen" / "Daniel F." / Binningen / Basel-Land)
# was a duplicate/erroneous entry; remove it. This shifts the remaining rows
# (old rows 9-16) up by one (new rows 8-15).
$ws.Rows.Item(8).Delete()

# The row that becomes the new row 8 ("à Wengen-Dörig" / Daniela / Bischofszell)
# had an incorrect "Bemerkung Selbstdispensation" flag; fix it from "Ja" to "Nein".
$ws.Cells.Item(8, 10).Value = "Nein"

# Append the new doctor/address/company row that was missing (row 16).
$ws.Cells.Item(16, 1).Value = "7601000010735"
$ws.Cells.Item(16, 2).Value = "Cevey"
$ws.Cells.Item(16, 3).Value = "Philippe Marc"
$ws.Cells.Item(16, 4).Value = "2000"
$ws.Cells.Item(16, 5).Value = "Neuchâtel"
$ws.Cells.Item(16, 6).Value = "Neuenburg"
$ws.Cells.Item(16, 7).Value = "CH"
$ws.Cells.Item(16, 8).Value = "Ärztin/Arzt"
$ws.Cells.Item(16, 9).Value = "Ja"
$ws.Cells.Item(16, 10).Value = "Nein"
